$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "browser" column is inserted as column C (pushing the former
# username/password/fname columns from C:E to D:F). Column widths are left
# untouched (not a real column insert), so cell values are rewritten
# directly instead of using EntireColumn.Insert().

# Row 1 (header)
$ws.Range("A1").Value = "testname"
$ws.Range("B1").Value = "execute"
$ws.Range("C1").Value = "browser"
$ws.Range("D1").Value = "username"
$ws.Range("E1").Value = "password"
$ws.Range("F1").Value = "fname"

# Row 2
$ws.Range("A2").Value = "loginLogoutTest"
$ws.Range("B2").Value = "yes"
$ws.Range("C2").Value = "chrome"
$ws.Range("D2").Value = "Admin"
$ws.Range("E2").Value = "admin123"
$ws.Range("F2").Value = "ScubaDrew615"

# Row 3
$ws.Range("A3").Value = "loginLogoutTest"
$ws.Range("B3").Value = "no"
$ws.Range("C3").Value = "chrome"
$ws.Range("D3").Value = "Admin"
$ws.Range("E3").Value = "admin123"
$ws.Range("F3").Value = "ScubaDrew615"

# Row 4
$ws.Range("A4").Value = "failedLoginLogoutTest"
$ws.Range("B4").Value = "yes"
$ws.Range("C4").Value = "chrome"
$ws.Range("D4").Value = "admin12"
$ws.Range("E4").Value = "admin123"
$ws.Range("F4").Value = "ScubaDrew615"

# Row 5
$ws.Range("A5").Value = "failedLoginLogoutTest"
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "chrome"
$ws.Range("D5").Value = "Admin"
$ws.Range("E5").Value = "Admin123"

# Row 6
$ws.Range("A6").Value = "loginLogoutTest"
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "chrome"
$ws.Range("D6").Value = "Admin"
$ws.Range("E6").Value = "admin123"

# Update selection to match the new active cell
$ws.Range("C5").Select()
